# Applies the timeline.docx edit described by the commit diff:
#  1. Insert a new "Landing pods (jack sees mountain range)" paragraph
#     (tab-indented) right before the "First probes" paragraph.
#  2. Move the <w:lastRenderedPageBreak/> marker from the "Screams"
#     paragraph to the "Jack runs into great creature..." paragraph.
#  3. Replace the "Jack befriends creature" paragraph with an expanded
#     block: a struck-through merge-conflict-style note about changing
#     the lion creature, the (now struck-through) "Jack befriends
#     creature" line, and several new "exploration of planet" notes.

$d = $word.ActiveDocument

function Get-PkgXml([string]$body) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $body + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Change 1: insert "Landing pods (jack sees mountain range)" paragraph
# before "First probes".
# ---------------------------------------------------------------------
$target = $d.Content.Find.Execute("First probes")
$firstProbesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "First probes`r") {
        $firstProbesPara = $d.Paragraphs($i)
        break
    }
}
$firstProbesPara.Range.InsertParagraphBefore()

$newPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "First probes`r") {
        $newPara = $d.Paragraphs($i - 1)
        break
    }
}
$body1 = '<w:body><w:p><w:r><w:tab/><w:t>Landing pods (jack sees mountain range)</w:t></w:r></w:p></w:body>'
$newPara.Range.InsertXML( (Get-PkgXml $body1) )

# ---------------------------------------------------------------------
# Change 2: move <w:lastRenderedPageBreak/> from "Screams" to
# "Jack runs into great creature, kills with a knife".
# ---------------------------------------------------------------------
$screamsPara = $null
$creaturePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq "Screams`r") { $screamsPara = $d.Paragraphs($i) }
    if ($t -eq "Jack runs into great creature, kills with a knife`r") { $creaturePara = $d.Paragraphs($i) }
}

$body2 = '<w:body><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Screams</w:t></w:r></w:p></w:body>'
$screamsPara.Range.InsertXML( (Get-PkgXml $body2) )

$body3 = '<w:body><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:lastRenderedPageBreak/><w:t>Jack runs into great creature, kills with a knife</w:t></w:r></w:p></w:body>'
$creaturePara.Range.InsertXML( (Get-PkgXml $body3) )

# ---------------------------------------------------------------------
# Change 3: replace "Jack befriends creature" paragraph with the
# expanded block of new paragraphs.
#
# The block ends with a bare empty <w:p/>. Word's Range.InsertXML drops
# a completely run-less trailing paragraph (it collapses into the
# replaced range's end boundary instead of staying distinct), so the
# replacement range is extended to also cover the following paragraph
# ("RRH finds jack...") and that paragraph's text is reproduced
# verbatim at the end of the inserted XML -- this keeps the empty
# paragraph intact as its own element.
# ---------------------------------------------------------------------
$befriendsPara = $null
$rrhPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Jack befriends creature`r") {
        $befriendsPara = $d.Paragraphs($i)
        $rrhPara = $d.Paragraphs($i + 1)
        break
    }
}

$combinedRange = $d.Range($befriendsPara.Range.Start, $rrhPara.Range.End)

$body4 = '<w:body>' +
    '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:strike/></w:rPr><w:t>&gt;&gt;&gt;&gt;&gt;</w:t></w:r>' +
        '<w:r><w:rPr><w:strike/></w:rPr><w:t>CHANGE LION CREATURE TO ALL ROCK</w:t></w:r>' +
        '<w:r><w:rPr><w:strike/></w:rPr><w:t>&lt;&lt;&lt;&lt;&lt;&lt;&lt;&lt;</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:strike/></w:rPr><w:t>Jack befriends creature</w:t></w:r>' +
    '</w:p>' +
    '<w:p/>' +
    '<w:p><w:r><w:t>Need some exploration of planet to fill space</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/><w:t>Ancient ruins?</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/><w:t>Previous civilization</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">Tau </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t>eti peoples (ancient type)??</w:t></w:r></w:p>' +
    '<w:p/>' +
    '<w:p><w:r><w:t>RRH finds jack, he releases creature upon them and slaughters them</w:t></w:r></w:p>' +
    '</w:body>'
$combinedRange.InsertXML( (Get-PkgXml $body4) )

Write-Host "Edits applied."
